$d = $word.ActiveDocument

# 1. Merge the "Javac" split runs (with spell-check proof errors) in the
#    Compilation paragraph into one continuous sentence. Word's Find/Replace
#    operates across run boundaries, so replacing the whole span collapses
#    it into a single run of text (dropping the proofErr spell-check marks).
$find1 = "Java source code is saved with the extension of ‘.java’ and it will pass through the Java compiler (Javac). Javac is a Java compiler that compiles Java code into a Java bytecode. Then, those Bytecode gets saved on the disk with the file extension ‘.class’. Java is an Object-Oriented Programming (OOP) language, therefore a program in Java is made of one or more classes. "
$replace1 = "Java source code is saved with the extension of ‘.java’ and it will pass through the Java compiler (Javac). Javac is a Java compiler that compiles Java code into a Java bytecode. Then, those Bytecode gets saved on the disk with the file extension ‘.class’. Java is an Object-Oriented Programming (OOP) language, therefore a program in Java is made of one or more classes. "

$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace1, 2)

# 2. Fix the typo "damages actions" -> "damage actions" in the Execution
#    paragraph.
$find2 = "performing damages actions"
$replace2 = "performing damage actions"

$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace2, 2)
